$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2901.8076
$ws.Range("J17").Value = 2901.8076
$ws.Range("L17").Value = 8705.4228
$ws.Range("N17").Value = -9041.4228

$ws.Range("H62").Value = 2112.0908
$ws.Range("I62").Value = 2267.5
$ws.Range("J62").Value = 1925.6
$ws.Range("K62").Value = 2267.5
$ws.Range("L62").Value = 1925.6
$ws.Range("M62").Value = -1643.5
$ws.Range("N62").Value = -3173.6

$ws.Range("H65").Value = 2112.0908
$ws.Range("I65").Value = 2267.5
$ws.Range("J65").Value = 1925.6
$ws.Range("K65").Value = 11337.5
$ws.Range("L65").Value = 9628
$ws.Range("M65").Value = -8217.5
$ws.Range("N65").Value = -15868

$ws.Range("H98").Value = 906
$ws.Range("I98").Value = 935.7895
$ws.Range("J98").Value = 811.6667
$ws.Range("K98").Value = 935.7895
$ws.Range("L98").Value = 811.6667
$ws.Range("M98").Value = 562.2105
$ws.Range("N98").Value = -3807.6667

$ws.Range("H106").Value = 3452.5
$ws.Range("I106").Value = 3405
$ws.Range("J106").Value = 3500
$ws.Range("K106").Value = 3405
$ws.Range("L106").Value = 3500
$ws.Range("M106").Value = -2774
$ws.Range("N106").Value = -4762

$ws.Range("H110").Value = 18940.4
$ws.Range("J110").Value = 18940.4
$ws.Range("L110").Value = 18940.4
$ws.Range("N110").Value = -27120.4

$ws.Range("H122").Value = 906
$ws.Range("I122").Value = 935.7895
$ws.Range("J122").Value = 811.6667
$ws.Range("K122").Value = 2807.3685
$ws.Range("L122").Value = 2435.0001
$ws.Range("M122").Value = -357.3685
$ws.Range("N122").Value = -7335.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 8642.666999999999
$ws.Range("I61").Value = 6243.273
$ws.Range("J61").Value = 15241
$ws.Range("K61").Value = 6243.273
$ws.Range("L61").Value = 15241
$ws.Range("M61").Value = -6031.273
$ws.Range("N61").Value = -15665

$ws.Range("H96").Value = 34333.332
$ws.Range("J96").Value = 34333.332
$ws.Range("L96").Value = 34333.332
$ws.Range("N96").Value = -39825.332

$ws.Range("H136").Value = 8642.666999999999
$ws.Range("I136").Value = 6243.273
$ws.Range("J136").Value = 15241
$ws.Range("K136").Value = 18729.819
$ws.Range("L136").Value = 45723
$ws.Range("M136").Value = -16179.819
$ws.Range("N136").Value = -50823

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 11691.375
$ws.Range("I82").Value = 2227.5
$ws.Range("J82").Value = 40083
$ws.Range("K82").Value = 2227.5
$ws.Range("L82").Value = 40083
$ws.Range("M82").Value = -1844.5
$ws.Range("N82").Value = -40849

$ws.Range("H85").Value = 11691.375
$ws.Range("I85").Value = 2227.5
$ws.Range("J85").Value = 40083
$ws.Range("K85").Value = 2227.5
$ws.Range("L85").Value = 40083
$ws.Range("M85").Value = -901.5
$ws.Range("N85").Value = -42735

$ws.Range("H134").Value = 28096.59
$ws.Range("I134").Value = 2207.25
$ws.Range("J134").Value = 93996.73
$ws.Range("K134").Value = 6621.75
$ws.Range("L134").Value = 281990.19
$ws.Range("M134").Value = -4086.75
$ws.Range("N134").Value = -287060.19

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5791.5806
$ws.Range("I31").Value = 6142.136
$ws.Range("J31").Value = 4934.6665
$ws.Range("K31").Value = 6142.136
$ws.Range("L31").Value = 4934.6665
$ws.Range("M31").Value = -5847.136
$ws.Range("N31").Value = -5524.6665

$ws.Range("H34").Value = 5791.5806
$ws.Range("I34").Value = 6142.136
$ws.Range("J34").Value = 4934.6665
$ws.Range("K34").Value = 6142.136
$ws.Range("L34").Value = 4934.6665
$ws.Range("M34").Value = -5940.136
$ws.Range("N34").Value = -5338.6665

$ws.Range("H58").Value = 1504775.1
$ws.Range("I58").Value = 2246137.5
$ws.Range("J58").Value = 3516.4
$ws.Range("K58").Value = 2246137.5
$ws.Range("L58").Value = 3516.4
$ws.Range("M58").Value = -2245934.5
$ws.Range("N58").Value = -3922.4

$ws.Range("H132").Value = 3247.621
$ws.Range("I132").Value = 3526.1162
$ws.Range("J132").Value = 2726.9565
$ws.Range("K132").Value = 10578.3486
$ws.Range("L132").Value = 8180.869499999999
$ws.Range("M132").Value = -8048.348599999999
$ws.Range("N132").Value = -13240.8695

$ws.Range("H136").Value = 1504775.1
$ws.Range("I136").Value = 2246137.5
$ws.Range("J136").Value = 3516.4
$ws.Range("K136").Value = 6738412.5
$ws.Range("L136").Value = 10549.2
$ws.Range("M136").Value = -6735862.5
$ws.Range("N136").Value = -15649.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 4905489.5
$ws.Range("I5").Value = 379.18182
$ws.Range("J5").Value = 13898192
$ws.Range("K5").Value = 1137.54546
$ws.Range("L5").Value = 41694576
$ws.Range("M5").Value = -1025.54546
$ws.Range("N5").Value = -41694800

$ws.Range("H18").Value = 9524180
$ws.Range("I18").Value = 10526620
$ws.Range("K18").Value = 31579860
$ws.Range("M18").Value = -31579691

$ws.Range("H68").Value = 5231.2915
$ws.Range("J68").Value = 6394.7896
$ws.Range("L68").Value = 19184.3688
$ws.Range("N68").Value = -20806.3688

$ws.Range("H71").Value = 5231.2915
$ws.Range("J71").Value = 6394.7896
$ws.Range("L71").Value = 57553.1064
$ws.Range("N71").Value = -65665.1064

$ws.Range("H112").Value = 2553
$ws.Range("J112").Value = 3970
$ws.Range("L112").Value = 11910
$ws.Range("N112").Value = -14126

$ws.Range("H122").Value = 742.25
$ws.Range("I122").Value = 524.9474
$ws.Range("J122").Value = 938.8570999999999
$ws.Range("K122").Value = 4724.5266
$ws.Range("L122").Value = 8449.713899999999
$ws.Range("M122").Value = -2274.5266
$ws.Range("N122").Value = -13349.7139

$ws.Range("H135").Value = 4905489.5
$ws.Range("I135").Value = 379.18182
$ws.Range("J135").Value = 13898192
$ws.Range("K135").Value = 3412.63638
$ws.Range("L135").Value = 125083728
$ws.Range("M135").Value = -877.6363799999999
$ws.Range("N135").Value = -125088798

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 210206400
$ws.Range("I10").Value = 262755000
$ws.Range("J10").Value = 12000
$ws.Range("K10").Value = 262755000
$ws.Range("L10").Value = 12000
$ws.Range("M10").Value = -262754831
$ws.Range("N10").Value = -12338

$ws.Range("H122").Value = 10377.857
$ws.Range("I122").Value = 14875.25
$ws.Range("J122").Value = 4381.3335
$ws.Range("K122").Value = 44625.75
$ws.Range("L122").Value = 13144.0005
$ws.Range("M122").Value = -42175.75
$ws.Range("N122").Value = -18044.0005

$ws.Range("H132").Value = 7108.025
$ws.Range("I132").Value = 5719.143
$ws.Range("J132").Value = 10348.75
$ws.Range("K132").Value = 17157.429
$ws.Range("L132").Value = 31046.25
$ws.Range("M132").Value = -14627.429
$ws.Range("N132").Value = -36106.25

$ws.Range("H136").Value = 19581
$ws.Range("J136").Value = 19581
$ws.Range("L136").Value = 58743
$ws.Range("N136").Value = -63843

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5000
$ws.Range("I40").Value = 5000
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 5000
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -4864
$ws.Range("N40").Value = -5272

$ws.Range("H82").Value = 1672.125
$ws.Range("I82").Value = 1466.6666
$ws.Range("J82").Value = 1795.4
$ws.Range("K82").Value = 1466.6666
$ws.Range("L82").Value = 1795.4
$ws.Range("M82").Value = -1105.6666
$ws.Range("N82").Value = -2517.4

$ws.Range("H85").Value = 1672.125
$ws.Range("I85").Value = 1466.6666
$ws.Range("J85").Value = 1795.4
$ws.Range("K85").Value = 1466.6666
$ws.Range("L85").Value = 1795.4
$ws.Range("M85").Value = -218.6666
$ws.Range("N85").Value = -4291.4

$ws.Range("H100").Value = 4610.6
$ws.Range("I100").Value = 2332.375
$ws.Range("J100").Value = 7214.2856
$ws.Range("K100").Value = 2332.375
$ws.Range("L100").Value = 7214.2856
$ws.Range("M100").Value = -1791.375
$ws.Range("N100").Value = -8296.285599999999

$ws.Range("H110").Value = 38000
$ws.Range("J110").Value = 38000
$ws.Range("L110").Value = 38000
$ws.Range("N110").Value = -46180

$ws.Range("H122").Value = 8133.0303
$ws.Range("I122").Value = 7938.4614
$ws.Range("J122").Value = 8855.714
$ws.Range("K122").Value = 23815.3842
$ws.Range("L122").Value = 26567.142
$ws.Range("M122").Value = -21365.3842
$ws.Range("N122").Value = -31467.142

$ws.Range("H130").Value = 279607.25
$ws.Range("J130").Value = 279607.25
$ws.Range("L130").Value = 279607.25
$ws.Range("N130").Value = -289647.25

$ws.Range("H132").Value = 6718.5186
$ws.Range("I132").Value = 7776.294
$ws.Range("J132").Value = 4920.3
$ws.Range("K132").Value = 23328.882
$ws.Range("L132").Value = 14760.9
$ws.Range("M132").Value = -20798.882
$ws.Range("N132").Value = -19820.9

$ws.Range("H136").Value = 4068.544
$ws.Range("I136").Value = 2520.5789
$ws.Range("J136").Value = 7164.4736
$ws.Range("K136").Value = 7561.736699999999
$ws.Range("L136").Value = 21493.4208
$ws.Range("M136").Value = -5011.736699999999
$ws.Range("N136").Value = -26593.4208

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 49250
$ws.Range("J101").Value = 49250
$ws.Range("L101").Value = 49250
$ws.Range("N101").Value = -55740

$ws.Range("H102").Value = 57000
$ws.Range("J102").Value = 57000
$ws.Range("L102").Value = 57000
$ws.Range("N102").Value = -63490

$ws.Range("H103").Value = 45301
$ws.Range("J103").Value = 45301
$ws.Range("L103").Value = 45301
$ws.Range("N103").Value = -47645

$ws.Range("H132").Value = 1993.5264
$ws.Range("I132").Value = 850.3333
$ws.Range("J132").Value = 3405.7058
$ws.Range("K132").Value = 2550.9999
$ws.Range("L132").Value = 10217.1174
$ws.Range("M132").Value = -20.9998999999998
$ws.Range("N132").Value = -15277.1174

$ws.Range("H136").Value = 4861.534
$ws.Range("I136").Value = 3647.4783
$ws.Range("J136").Value = 6929.926
$ws.Range("K136").Value = 10942.4349
$ws.Range("L136").Value = 20789.778
$ws.Range("M136").Value = -8392.4349
$ws.Range("N136").Value = -25889.778
